$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 new columns before column D (existing D..J shift to P..V)
$ws.Range("D1:O1").EntireColumn.Insert()

# Fill in the new header cells (row 1) with the "Unnamed: 0.1.1..." series
$headers = @(
    "Unnamed: 0.1.1",
    "Unnamed: 0.1.1.1",
    "Unnamed: 0.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1",
    "Unnamed: 0.1.1.1.1.1.1.1.1.1.1.1.1.1"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 4 + $i).Value2 = $headers[$i]
}

# Fill in the new data cells (rows 2 and 3) matching column C's value
$ws.Range("D2:O2").Value2 = $ws.Range("C2").Value2
$ws.Range("D3:O3").Value2 = $ws.Range("C3").Value2
